# Apply updated cryptocurrency price/volume data to Sheet1.
# Cell values are written via Formula with a leading apostrophe so that
# numeric-looking strings (prices, percentages) are preserved as text,
# matching the inline-string storage used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'94.209.37"
$ws.Range("E2").Formula = "'  +1.63%  "
$ws.Range("D3").Formula = "'3.082.57"
$ws.Range("E3").Formula = "'  -1.10%  "
$ws.Range("E4").Formula = "'  +0.01%  "
$ws.Range("D5").Formula = "'233.76"
$ws.Range("E5").Formula = "'  -3.40%  "
$ws.Range("D6").Formula = "'607.96"
$ws.Range("E6").Formula = "'  -1.23%  "
$ws.Range("D7").Formula = "'1.09"
$ws.Range("E7").Formula = "'  -0.22%  "
$ws.Range("D8").Formula = "'0.376"
$ws.Range("E8").Formula = "'  -6.08%  "
$ws.Range("E9").Formula = "'  -0.05%  "
$ws.Range("D10").Formula = "'0.817"
$ws.Range("E10").Formula = "'  +11.89%  "
$ws.Range("D11").Formula = "'3.081.35"
$ws.Range("E11").Formula = "'  -1.05%  "
$ws.Range("E12").Formula = "'  -3.22%  "
$ws.Range("D13").Formula = "'94.104.66"
$ws.Range("E13").Formula = "'  +1.62%  "
$ws.Range("E14").Formula = "'  -5.62%  "
$ws.Range("D15").Formula = "'33.84"
$ws.Range("E15").Formula = "'  -1.62%  "
$ws.Range("D16").Formula = "'5.26"
$ws.Range("E16").Formula = "'  -4.20%  "
$ws.Range("D17").Formula = "'3.651.92"
$ws.Range("E17").Formula = "'  -1.37%  "
$ws.Range("D18").Formula = "'3.049.52"
$ws.Range("E18").Formula = "'  -2.27%  "
$ws.Range("D19").Formula = "'3.59"
$ws.Range("E19").Formula = "'  -3.89%  "
$ws.Range("D20").Formula = "'14.48"
$ws.Range("E20").Formula = "'  -1.36%  "
$ws.Range("D21").Formula = "'5.71"
$ws.Range("E21").Formula = "'  -1.26%  "
$ws.Range("D22").Formula = "'439.53"
$ws.Range("E22").Formula = "'  -1.85%  "
$ws.Range("D23").Formula = "'8.78"
$ws.Range("E23").Formula = "'  -6.94%  "
$ws.Range("D24").Formula = "'0.0000190"
$ws.Range("E24").Formula = "'  -7.55%  "
$ws.Range("E25").Formula = "'  +5.05%  "
$ws.Range("D26").Formula = "'5.52"
$ws.Range("E26").Formula = "'  -4.71%  "
$ws.Range("D27").Formula = "'84.96"
$ws.Range("E27").Formula = "'  -2.20%  "
$ws.Range("D28").Formula = "'11.85"
$ws.Range("E28").Formula = "'  +0.82%  "
$ws.Range("D29").Formula = "'3.253.35"
$ws.Range("E29").Formula = "'  -0.93%  "
$ws.Range("D30").Formula = "'0.999"
$ws.Range("E30").Formula = "'  -0.37%  "
$ws.Range("D31").Formula = "'0.246"
$ws.Range("E31").Formula = "'  +6.59%  "
$ws.Range("D32").Formula = "'0.177"
$ws.Range("E32").Formula = "'  +4.97%  "
$ws.Range("E33").Formula = "'  -10.60%  "
$ws.Range("D34").Formula = "'9.10"
$ws.Range("E34").Formula = "'  -1.68%  "
$ws.Range("E35").Formula = "'  -0.68%  "
$ws.Range("D36").Formula = "'7.71"
$ws.Range("E36").Formula = "'  -3.54%  "
$ws.Range("E37").Formula = "'  -4.33%  "
$ws.Range("D38").Formula = "'25.54"
$ws.Range("E38").Formula = "'  -2.55%  "
$ws.Range("E39").Formula = "'  -2.05%  "
$ws.Range("D40").Formula = "'0.442"
$ws.Range("E40").Formula = "'  +0.84%  "
$ws.Range("D41").Formula = "'23.99"
$ws.Range("E41").Formula = "'  +4.00%  "
$ws.Range("B42").Formula = "'Bittensor"
$ws.Range("C42").Formula = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Formula = "'465.86"
$ws.Range("E42").Formula = "'  -3.68%  "
$ws.Range("B43").Formula = "'MantraDAO"
$ws.Range("C43").Formula = "'https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").Formula = "'3.71"
$ws.Range("E43").Formula = "'  -11.27%  "
$ws.Range("E44").Formula = "'  -4.36%  "
$ws.Range("D46").Formula = "'3.10"
$ws.Range("E46").Formula = "'  -11.54%  "
$ws.Range("D47").Formula = "'159.78"
$ws.Range("E47").Formula = "'  -1.02%  "
$ws.Range("D48").Formula = "'1.84"
$ws.Range("E48").Formula = "'  -4.34%  "
$ws.Range("D49").Formula = "'0.675"
$ws.Range("D50").Formula = "'43.78"
$ws.Range("E50").Formula = "'  -0.74%  "
$ws.Range("D51").Formula = "'0.998"
$ws.Range("E51").Formula = "'  +0.01%  "
